# Progress and time spent update for MasterDocumentSummaryV1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 43 (DespatchingVx / Outbound scheduling & control) ---
# Budget 80% -> 95%, Time used 14 -> 18 hrs
$ws.Cells.Item(43, 4).Value2 = 0.95
$ws.Cells.Item(43, 8).Value2 = 18

# --- Row 44 (DoorScanningVx / Control freight at courier handover) ---
# Budget 80% -> 95%, Time used 9 -> 12 hrs
$ws.Cells.Item(44, 4).Value2 = 0.95
$ws.Cells.Item(44, 8).Value2 = 12

# --- Row 46 (ClaimsVx / Claims on supplier tracking and payment release) ---
# Progress "To Start" cleared, Time Needed 20 -> 2, Time used 0 -> 2
# Review note updated & now wraps
$ws.Cells.Item(46, 4).Clear()
$ws.Cells.Item(46, 5).Value2 = 2
$ws.Cells.Item(46, 8).Value2 = 2
$ws.Cells.Item(46, 10).Value2 = "Reduced time from 20 hours - completed basic overview narrative`nNon trivial solution due to supplier settlement permutations"
$ws.Cells.Item(46, 10).WrapText = $true

# --- Row 47 (PODScanningVx / Scan capture signed POD's) ---
# Progress "To Start" cleared, Time Needed 14 -> 2, Time used 0 -> 2
# New review note added
$ws.Cells.Item(47, 4).Clear()
$ws.Cells.Item(47, 5).Value2 = 2
$ws.Cells.Item(47, 8).Value2 = 2
$ws.Cells.Item(47, 10).Value2 = "Time reduced from 14 hours. Limited narative provided"

# --- Sheet view: selection moved to J48 ---
$ws.Activate()
$ws.Range("J48").Select()
